# Applies the cryptos.xlsx price/volume updates described in the commit
# "Updated cryptos list on Fri May 26 15:30:07 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.989.63"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "1.848.56"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("E4").Value = "  +0.11%  "
$cellStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.50"
$ws.Range("D5").Style = $cellStyle
$ws.Range("E5").Value = "  +1.22%  "
$cellStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = $cellStyle
$ws.Range("E6").Value = "  +0.10%  "
$cellStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4676"
$ws.Range("D7").Style = $cellStyle
$ws.Range("E7").Value = "  +3.26%  "
$cellStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3637"
$ws.Range("D8").Style = $cellStyle
$ws.Range("E8").Value = "  +0.97%  "
$cellStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07178"
$ws.Range("D9").Style = $cellStyle
$ws.Range("E9").Value = "  +1.42%  "
$cellStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9308"
$ws.Range("D10").Style = $cellStyle
$ws.Range("E10").Value = "  +4.16%  "
$cellStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.61"
$ws.Range("D11").Style = $cellStyle
$ws.Range("E11").Value = "  +0.83%  "
$cellStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07680"
$ws.Range("D12").Style = $cellStyle
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("D13").Value = "1.825.17"
$ws.Range("E13").Value = "  +1.01%  "
$cellStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.297"
$ws.Range("D14").Style = $cellStyle
$cellStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.418"
$ws.Range("D15").Style = $cellStyle
$ws.Range("E15").Value = "  +1.54%  "
$cellStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.24"
$ws.Range("D16").Style = $cellStyle
$ws.Range("E16").Value = "  +2.98%  "
$cellStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("D17").Style = $cellStyle
$ws.Range("E17").Value = "  +0.23%  "
$cellStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008594"
$ws.Range("D18").Style = $cellStyle
$ws.Range("E18").Value = "  +1.16%  "
$cellStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("D19").Style = $cellStyle
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "27.013.15"
$cellStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.41"
$ws.Range("D21").Style = $cellStyle
$ws.Range("E21").Value = "  +1.45%  "
$cellStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.038"
$ws.Range("D22").Style = $cellStyle
$ws.Range("E22").Value = "  +1.36%  "
$ws.Range("E23").Value = "  +1.16%  "
$cellStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.940"
$ws.Range("D24").Style = $cellStyle
$ws.Range("E24").Value = "  -0.93%  "
$cellStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.65"
$ws.Range("D25").Style = $cellStyle
$ws.Range("E25").Value = "  +0.30%  "
$cellStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.05"
$ws.Range("D26").Style = $cellStyle
$ws.Range("E26").Value = "  +1.39%  "
$cellStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.035"
$ws.Range("D27").Style = $cellStyle
$ws.Range("E27").Value = "  -1.63%  "
$cellStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.09"
$ws.Range("D28").Style = $cellStyle
$ws.Range("E28").Value = "  +1.78%  "
$cellStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.932"
$ws.Range("D29").Style = $cellStyle
$ws.Range("E29").Value = "  +1.53%  "
$cellStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08867"
$ws.Range("D30").Style = $cellStyle
$ws.Range("E30").Value = "  +1.94%  "
$cellStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.183"
$ws.Range("D31").Style = $cellStyle
$ws.Range("E31").Value = "  +2.07%  "
$cellStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.859"
$ws.Range("D32").Style = $cellStyle
$ws.Range("E32").Value = "  +0.66%  "
$cellStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.180"
$ws.Range("D33").Style = $cellStyle
$ws.Range("E33").Value = "  +6.71%  "
$cellStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7456"
$ws.Range("D34").Style = $cellStyle
$ws.Range("E34").Value = "  +2.90%  "
$cellStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.475"
$ws.Range("D35").Style = $cellStyle
$ws.Range("E35").Value = "  +0.17%  "
$cellStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.088"
$ws.Range("D36").Style = $cellStyle
$ws.Range("E36").Value = "  +1.28%  "
$cellStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.995"
$ws.Range("D37").Style = $cellStyle
$ws.Range("E37").Value = "  +2.95%  "
$cellStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01942"
$ws.Range("D38").Style = $cellStyle
$cellStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05174"
$ws.Range("D39").Style = $cellStyle
$ws.Range("E39").Value = "  +1.27%  "
$cellStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5139"
$ws.Range("D40").Style = $cellStyle
$ws.Range("E40").Value = "  +0.42%  "
$cellStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.901"
$ws.Range("D41").Style = $cellStyle
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("E42").Value = "  -0.13%  "
$cellStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.179"
$ws.Range("D43").Style = $cellStyle
$ws.Range("E43").Value = "  +1.88%  "
$cellStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.49"
$ws.Range("D44").Style = $cellStyle
$ws.Range("E44").Value = "  +5.15%  "
$cellStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4714"
$ws.Range("D45").Style = $cellStyle
$ws.Range("E45").Value = "  +0.84%  "
$cellStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.009"
$ws.Range("D46").Style = $cellStyle
$ws.Range("E46").Value = "  +0.22%  "
$cellStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.26"
$ws.Range("D47").Style = $cellStyle
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("E49").Value = "  +1.07%  "
$cellStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.54"
$ws.Range("D50").Style = $cellStyle
$ws.Range("E50").Value = "  +1.36%  "
$cellStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.14"
$ws.Range("D51").Style = $cellStyle
$ws.Range("E51").Value = "  +0.00%  "

Write-Output "Applied cryptos list update (92 cells changed)"
